$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# Row 7 ("Experimental") previously had an empty Value cell (B7).
# Set it to the literal text "false" (not the boolean FALSE).
# A bare Value assignment of "false" gets auto-coerced to a boolean by the
# engine (same as typing it into Excel), so force text entry with a leading
# apostrophe, then re-apply B7's original formatting (copied from the
# neighbouring A7 cell, which still carries the unmodified style) so the
# quote-prefix marker doesn't leave a stray style behind.
$meta.Range("B7").Value = "'false"
$meta.Range("A7").Copy()
$meta.Range("B7").PasteSpecial(-4122)  # xlPasteFormats

# Row 8 ("Date"): update the generation timestamp.
$meta.Range("B8").Value = "2025-11-30T13:08:37+00:00"
